$d = $word.ActiveDocument

# Merge paragraph runs that were split across multiple <w:r> elements
# back into a single run per paragraph (content unchanged, only run
# boundaries collapse). Find/Replace with identical old/new text that
# spans a run boundary causes the runtime to merge the touched runs.
$mergedTexts = @(
    "Lorem ipsum dolor sit amet, consectetur adipiscing elit. Vivamus fermentum id lacus vel vulputate. Cras malesuada dolor sit amet orci mattis placerat. Ut eget aliquam leo, at consequat velit. Donec convallis est quis molestie sollicitudin. Sed sodales libero leo, nec pretium neque venenatis eu. Sed a nibh pharetra purus vehicula eleifend eu sit amet velit. Nulla ut neque eget tortor placerat luctus. Proin eu tortor quis mi iaculis rhoncus.",
    "Morbi lacinia, magna et sollicitudin fermentum, odio nisi consectetur erat, ac luctus erat elit eu enim. Phasellus tortor nulla, egestas ac dui sit amet, dictum viverra enim. Maecenas venenatis ligula et finibus ornare. Vivamus fermentum arcu a ipsum mattis rhoncus. Aliquam pellentesque nisi tristique faucibus finibus. Donec non ipsum felis. In ultricies mollis augue vel euismod. Vivamus varius, tortor non egestas aliquam, est lectus laoreet nunc, in varius felis elit ac ex. Mauris nec mauris at ipsum elementum laoreet. Aenean vitae libero ex. Donec rhoncus lacinia libero ac consectetur :",
    "Vestibulum sagittis libero in odio dignissim, tincidunt varius lacus finibus. Nullam ut vehicula magna. Nam nunc urna, ultricies ut nunc sed, accumsan molestie enim. Duis non augue pharetra, tristique nisi at, feugiat tellus. Nulla ultricies sapien quis aliquet mattis. Nunc elementum nec metus ac ornare. Fusce a diam ante. Maecenas id massa et purus semper iaculis a vitae orci. Integer euismod ipsum nec quam semper consequat. Class aptent taciti sociosqu ad litora torquent per conubia nostra, per inceptos himenaeos. Quisque a lorem lorem. Fusce nec iaculis lectus. Cras nulla dolor, ultrices at odio in, fringilla efficitur nunc.",
    "Suspendisse neque turpis, facilisis ut efficitur vulputate, tempor pellentesque diam. Nam viverra viverra odio quis aliquam. Duis lobortis ipsum libero, sit amet gravida tortor rhoncus a. In lobortis eu nibh vitae vehicula. Phasellus lacus tortor, venenatis ac nunc at, faucibus tempor urna. Donec lobortis bibendum nibh sit amet commodo. Etiam gravida mauris et enim efficitur ultricies. Lorem ipsum dolor sit amet, consectetur adipiscing elit. Proin lobortis ultrices urna, at consectetur ipsum pharetra eu.",
    "In sit amet mattis purus. Vivamus ornare pharetra lacus at consequat. In ullamcorper, magna sagittis dignissim egestas, libero nibh ornare nunc, ac commodo sem sem a odio. Sed at nulla sit amet purus egestas auctor eu fringilla nulla. Vivamus eu convallis dolor. Mauris aliquam nec tellus a dignissim. Quisque pharetra, tortor eu tincidunt cursus, nibh dolor vehicula odio, vel fermentum ligula felis nec ligula. Vestibulum vel odio eget elit feugiat convallis. Nulla non sem ultrices, pretium ligula nec, aliquam mi. Integer blandit nisl nec enim facilisis blandit. Nullam feugiat risus vitae nibh laoreet, sit amet condimentum ipsum auctor. Nunc ut dui vel nisl malesuada consectetur at vitae mi. Duis nec leo at orci bibendum porttitor. Maecenas accumsan ac est id vehicula. Mauris non nisl nisl.",
    "Phasellus leo dolor, auctor non scelerisque vitae, convallis id enim. Vestibulum nec quam et velit sodales pellentesque. Vivamus sodales interdum iaculis. Fusce sed orci tellus. Etiam vel sagittis felis. Nulla eu aliquet ligula. Nunc posuere est eget turpis condimentum pellentesque. Curabitur arcu nulla, auctor quis nibh quis, gravida blandit quam. Vivamus facilisis nunc mi, scelerisque pellentesque lacus volutpat quis.",
    "Cras at enim orci. Curabitur interdum rhoncus sem eget cursus. Pellentesque convallis tincidunt odio, eu lobortis lacus blandit et. Maecenas eros ex, volutpat malesuada sem et, luctus bibendum felis. Duis non facilisis neque. In cursus est a blandit tristique. Phasellus sed vehicula metus. Morbi quis posuere augue, quis iaculis enim. Integer consequat lacus vel facilisis bibendum. Sed nunc lectus, aliquet ut pellentesque eget, mollis in ligula. Donec mollis cursus rhoncus.",
    "Aenean eget ullamcorper leo. In varius posuere sapien et efficitur. Pellentesque maximus elit ut iaculis porttitor. Integer fermentum urna tristique lacus lacinia feugiat. Aliquam eu augue leo. Cras venenatis pellentesque nibh, imperdiet fermentum libero congue a. Nullam ac ornare dolor. Curabitur turpis orci, maximus vehicula suscipit a, maximus sit amet tortor. Fusce elementum nisl vel lectus lacinia, sed lobortis odio luctus. Sed pretium porta sem ac suscipit. Morbi nec nunc sed nisi interdum scelerisque. Donec ultrices leo in arcu gravida volutpat. Mauris sed blandit metus, in egestas libero. Ut ac lacinia magna. Maecenas volutpat, velit a dignissim rhoncus, velit lorem euismod dui, vel vehicula orci augue eget odio. Nullam nibh turpis, viverra eget luctus ac, semper ac est.",
    "Nullam ac volutpat felis. Nullam consectetur enim nisl, non pharetra neque imperdiet non. Donec non dictum risus. Integer rutrum condimentum commodo. Integer maximus neque non eros sollicitudin, in sodales nisi malesuada. Praesent pellentesque dapibus viverra. In in vehicula elit, in luctus nisi.",
)

foreach ($t in $mergedTexts) {
    $null = $d.Content.Find.Execute($t, $true, $false, $false, $false, $false, $true, 1, $false, $t, 2)
}

# Paragraph with <w:lastRenderedPageBreak/> mid-paragraph: merge the three
# runs preceding the break normally, then merge the run *after* the break
# by deleting its text and re-inserting it with Range.InsertAfter, which
# appends a new trailing run without disturbing the run that owns the
# page-break marker (a plain cross-run Find/Replace would merge that
# owning run away and silently drop the marker).
$beforeBreak = "Aliquam sem est, faucibus non est eget, laoreet mollis dolor. Nulla in eros felis. Cras id est pretium, fermentum tortor et, luctus sapien. Maecenas eget mattis quam. Vestibulum urna tortor, egestas quis nisl nec, convallis suscipit turpis. Phasellus pretium, nulla semper lacinia faucibus, nisi lorem bibendum nisi, id rhoncus elit lacus vel neque. Fusce finibus dolor sit amet aliquam facilisis. Fusce imperdiet purus vitae elit suscipit, ut tempus nunc "
$null = $d.Content.Find.Execute($beforeBreak, $true, $false, $false, $false, $false, $true, 1, $false, $beforeBreak, 2)

$afterTail = "met consectetur tellus egestas nec."
$tailRange = $d.Content.Duplicate
$null = $tailRange.Find.Execute($afterTail, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$tailRange.Text = ""
$pageBreakPara = $tailRange.Paragraphs.Item(1)
$pageBreakPara.Range.InsertAfter($afterTail)
